$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header label B1 from [ug] to [mg]
$ws.Range("B1").Value = "Tryptophane par semaine [mg]"

# New Tryptophan weekly values (mg) for rows 2-79 (CID 2-79)
$values = @(
    6881.415, 8018.889999999999, 4335.264999999999, 3046.115, 6858.355, 5483.139999999999, 8673.16, 6679.349999999999, 8598.4, 10306.6,
    8123.799999999999, 5425.799999999999, 5232.24, 4175.955, 8011.215, 9916.864999999998, 9873.014999999999, 14342.67, 6737.09, 2925.455,
    7012.285000000001, 3738.14, 8070.21, 5794.94, 4977.740000000001, 4071.925, 5702.205000000001, 7215.03, 12451.665, 7949.439999999999,
    3834.19, 8814.665000000001, 4433.174999999999, 5652.44, 3859.74, 11352.54, 19785.1, 7133.440000000001, 2524.675, 6752.215,
    5276.48, 6943.8, 6810.45, 8086.34, 6737.465, 5889.865, 2820.975, 8982.115, 3532.065, 4052.98,
    5017.425, 3643.845, 10420.415, 7751.515, 4542.865000000001, 6137.29, 11813.315, 4627.865, 3811.64, 7972.184999999999,
    6158.135, 2861.805, 2668.665, 7557.79, 4012.06, 6887.025, 6466.515, 3394.315, 2913.825, 6564.974999999999,
    6001.3, 2885.94, 6324.639999999999, 6347.974999999999, 6278.205, 8396.58, 5322.1, 7380.800000000001
)

for ($i = 0; $i -lt $values.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $values[$i]
}
